$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - update the normalized-length bin values for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data updates
$ws.Range("B2").Value = 92.993913932287896
$ws.Range("C2").Value = 93.59849480056269
$ws.Range("D2").Value = 94.079494538077839
$ws.Range("E2").Value = 94.133727192417666

# Row 3 data updates
$ws.Range("B3").Value = 94.067146827857201
$ws.Range("C3").Value = 96.781547447244463
$ws.Range("D3").Value = 94.798422624694567
$ws.Range("E3").Value = 94.264979334286821

# Update the selection to match the new range B1:E3
$ws.Range("B1:E3").Select()
